$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 28.333334
$ws.Range("I6").Value = 28.333334
$ws.Range("K6").Value = 85.00000199999999
$ws.Range("M6").Value = 26.99999800000001
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("N28").ClearContents()
$ws.Range("H33").Value = 147
$ws.Range("I33").Value = 125
$ws.Range("J33").Value = 185.5
$ws.Range("K33").Value = 125
$ws.Range("L33").Value = 185.5
$ws.Range("M33").Value = 104
$ws.Range("N33").Value = -643.5
$ws.Range("H43").Value = 2849.5
$ws.Range("H113").Value = 3499.6667
$ws.Range("I113").Value = 2832.6667
$ws.Range("J113").Value = 4166.6665
$ws.Range("K113").Value = 2832.6667
$ws.Range("L113").Value = 4166.6665
$ws.Range("M113").Value = 421.3332999999998
$ws.Range("N113").Value = -10674.6665
$ws.Range("H116").Value = 7979.8
$ws.Range("I116").Value = 8000
$ws.Range("J116").Value = 7974.75
$ws.Range("K116").Value = 8000
$ws.Range("L116").Value = 7974.75
$ws.Range("M116").Value = -4558
$ws.Range("N116").Value = -14858.75
$ws.Range("H137").Value = 1817.0769
$ws.Range("I137").Value = 1524.7273
$ws.Range("J137").Value = 3425
$ws.Range("K137").Value = 4574.1819
$ws.Range("L137").Value = 10275
$ws.Range("M137").Value = -2024.1819
$ws.Range("N137").Value = -15375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 671.875
$ws.Range("I2").Value = 435
$ws.Range("J2").Value = 1698.3334
$ws.Range("K2").Value = 435
$ws.Range("L2").Value = 1698.3334
$ws.Range("M2").Value = -322
$ws.Range("N2").Value = -1924.3334
$ws.Range("H45").Value = 6649.9
$ws.Range("I45").Value = 6649.9
$ws.Range("K45").Value = 6649.9
$ws.Range("M45").Value = -6272.9
$ws.Range("H62").Value = 53275.668
$ws.Range("H65").Value = 53275.668
$ws.Range("H116").Value = 671.875
$ws.Range("I116").Value = 435
$ws.Range("J116").Value = 1698.3334
$ws.Range("K116").Value = 435
$ws.Range("L116").Value = 1698.3334
$ws.Range("M116").Value = 1859
$ws.Range("N116").Value = -6286.3334
$ws.Range("H122").Value = 2633.3333
$ws.Range("I122").Value = 1700
$ws.Range("K122").Value = 5100
$ws.Range("M122").Value = -2650

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 671.875
$ws.Range("I3").Value = 435
$ws.Range("J3").Value = 1698.3334
$ws.Range("K3").Value = 435
$ws.Range("L3").Value = 1698.3334
$ws.Range("M3").Value = -321
$ws.Range("N3").Value = -1926.3334
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5147.25
$ws.Range("I99").Value = 4464.8
$ws.Range("K99").Value = 4464.8
$ws.Range("M99").Value = -2966.8
$ws.Range("H122").Value = 2784
$ws.Range("I122").Value = 2784
$ws.Range("K122").Value = 8352
$ws.Range("M122").Value = -5902
$ws.Range("H126").Value = 5147.25
$ws.Range("I126").Value = 4464.8
$ws.Range("K126").Value = 13394.4
$ws.Range("M126").Value = -10924.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2332.5
$ws.Range("I5").Value = 2332.5
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 6997.5
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -6885.5
$ws.Range("N5").ClearContents()
$ws.Range("H6").Value = 350
$ws.Range("I6").Value = 350
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1050
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -937
$ws.Range("N6").ClearContents()
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()
$ws.Range("H37").Value = 59996
$ws.Range("J37").Value = 59996
$ws.Range("L37").Value = 179988
$ws.Range("N37").Value = -180212
$ws.Range("H121").Value = 17506.1
$ws.Range("J121").Value = 6878.875
$ws.Range("L121").Value = 20636.625
$ws.Range("N121").Value = -23256.625
$ws.Range("H131").Value = 325539.53
$ws.Range("I131").Value = 2333
$ws.Range("J131").Value = 360168.78
$ws.Range("K131").Value = 6999
$ws.Range("L131").Value = 1080506.34
$ws.Range("M131").Value = -1959
$ws.Range("N131").Value = -1090586.34
$ws.Range("H135").Value = 2332.5
$ws.Range("I135").Value = 2332.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 20992.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -18457.5
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 4063
$ws.Range("I136").Value = 4063
$ws.Range("K136").Value = 12189
$ws.Range("M136").Value = -7089
$ws.Range("H141").Value = 13124
$ws.Range("I141").Value = 13124
$ws.Range("K141").Value = 39372
$ws.Range("M141").Value = -34192

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7111.1035
$ws.Range("I70").Value = 6900.5557
$ws.Range("J70").Value = 7455.636
$ws.Range("K70").Value = 6900.5557
$ws.Range("L70").Value = 7455.636
$ws.Range("M70").Value = -6630.5557
$ws.Range("N70").Value = -7995.636
$ws.Range("H73").Value = 7111.1035
$ws.Range("I73").Value = 6900.5557
$ws.Range("J73").Value = 7455.636
$ws.Range("K73").Value = 6900.5557
$ws.Range("L73").Value = 7455.636
$ws.Range("M73").Value = -5964.5557
$ws.Range("N73").Value = -9327.636
$ws.Range("H102").Value = 1649.5
$ws.Range("I102").Value = 1757.2
$ws.Range("J102").Value = 1111
$ws.Range("K102").Value = 1757.2
$ws.Range("L102").Value = 1111
$ws.Range("M102").Value = -135.2
$ws.Range("N102").Value = -4355
$ws.Range("H113").Value = 627.5
$ws.Range("I113").Value = 475
$ws.Range("J113").Value = 780
$ws.Range("K113").Value = 475
$ws.Range("L113").Value = 780
$ws.Range("M113").Value = 1695
$ws.Range("N113").Value = -5120
$ws.Range("H122").Value = 1513.7142
$ws.Range("I122").Value = 1779.6
$ws.Range("J122").Value = 849
$ws.Range("K122").Value = 5338.799999999999
$ws.Range("L122").Value = 2547
$ws.Range("M122").Value = -2888.799999999999
$ws.Range("N122").Value = -7447

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6016.2964
$ws.Range("I7").Value = 3419.2307
$ws.Range("J7").Value = 8427.857
$ws.Range("K7").Value = 3419.2307
$ws.Range("L7").Value = 8427.857
$ws.Range("M7").Value = -3307.2307
$ws.Range("N7").Value = -8651.857
$ws.Range("H22").Value = 4380.2
$ws.Range("I22").Value = 4967
$ws.Range("J22").Value = 3500
$ws.Range("K22").Value = 4967
$ws.Range("L22").Value = 3500
$ws.Range("M22").Value = -4672
$ws.Range("N22").Value = -4090
$ws.Range("H27").Value = 4380.2
$ws.Range("I27").Value = 4967
$ws.Range("J27").Value = 3500
$ws.Range("K27").Value = 4967
$ws.Range("L27").Value = 3500
$ws.Range("M27").Value = -4860
$ws.Range("N27").Value = -3714
$ws.Range("H55").Value = 137.38461
$ws.Range("I55").Value = 78.59999999999999
$ws.Range("K55").Value = 78.59999999999999
$ws.Range("M55").Value = 94.40000000000001
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H122").Value = 5233.5127
$ws.Range("I122").Value = 4354.615
$ws.Range("J122").Value = 6991.3076
$ws.Range("K122").Value = 13063.845
$ws.Range("L122").Value = 20973.9228
$ws.Range("M122").Value = -10613.845
$ws.Range("N122").Value = -25873.9228
$ws.Range("H126").Value = 6016.2964
$ws.Range("I126").Value = 3419.2307
$ws.Range("J126").Value = 8427.857
$ws.Range("K126").Value = 10257.6921
$ws.Range("L126").Value = 25283.571
$ws.Range("M126").Value = -7787.6921
$ws.Range("N126").Value = -30223.571
$ws.Range("H132").Value = 2831.75
$ws.Range("I132").Value = 2497.3333
$ws.Range("K132").Value = 7491.999899999999
$ws.Range("M132").Value = -4961.999899999999
$ws.Range("H136").Value = 2970.4443
$ws.Range("I136").Value = 1955.8334
$ws.Range("K136").Value = 5867.5002
$ws.Range("M136").Value = -3317.5002

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 745
$ws.Range("I122").Value = 745
$ws.Range("K122").Value = 2235
$ws.Range("M122").Value = 215
